$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.427.42'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -2.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.834.60'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.68%  '
$ws.Range('E4').Value = '  -0.92%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.70'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4609'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -3.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3814'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.31'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07911'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9737'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -4.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.10'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.79%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.901'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.62%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.818.15'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.031'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06610'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001028'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.03'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.443.76'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.356'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.86'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.36%  '
$ws.Range('E25').Value = '  -1.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.18'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.40'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.066'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.328'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.86'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.9535'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09290'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.01%  '
$ws.Range('E33').Value = '  -1.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.247'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.317'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05938'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02191'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.078'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.159'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5795'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.24%  '
$ws.Range('E41').Value = '  -2.70%  '
$ws.Range('E42').Value = '  -3.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.261'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5485'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.33%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.97'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.866'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06656'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '109.89'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.42%  '
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000289'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.87%  '
